# 19 - Formatação condicional - regras de realce
# Applies:
#  - a few data corrections (B2, B4, C2, C8)
#  - conditional formatting "highlight cell rules" on B2:B17 (Nacional/Internacional)
#    and C2:C17 (stock level thresholds)
#  - restores the selection to C9

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data corrections -------------------------------------------------
$ws.Range("B2").Value = "Nacional"
$ws.Range("B4").Value = "Internacional"
$ws.Range("C2").Value = 29000
$ws.Range("C8").Value = 120000

# The "Origem" column (B) carried a leftover no-op fill style; clear it so
# the column falls back to the default/Normal style before layering the
# new conditional formatting rules on top.
$ws.Range("B2:B17").Style = "Normal"

# --- Conditional formatting: Coluna de Origem (B2:B17) -----------------
$origem = $ws.Range("B2:B17")

# "Yellow Fill with Dark Yellow Text" - Nacional
$fcNacional = $origem.FormatConditions.Add(1, 3, '"Nacional"')
$fcNacional.Font.Color = 22428
$fcNacional.Interior.Color = 10284031

# "Light Red Fill with Dark Red Text" - Internacional
$fcInternacional = $origem.FormatConditions.Add(1, 3, '"Internacional"')
$fcInternacional.Font.Color = 393372
$fcInternacional.Interior.Color = 13551615

# --- Conditional formatting: Quantidade em Estoque (C2:C17) ------------
$estoque = $ws.Range("C2:C17")

# "Green Fill with Dark Green Text" - Acima de 120.000 (Saudável)
$fcAlto = $estoque.FormatConditions.Add(1, 5, 120000)
$fcAlto.Font.Color = 24832
$fcAlto.Interior.Color = 13561798

# "Light Red Fill with Dark Red Text" - Abaixo de 30.000 (Crítico)
$fcBaixo = $estoque.FormatConditions.Add(1, 6, 30000)
$fcBaixo.Font.Color = 393372
$fcBaixo.Interior.Color = 13551615

# "Yellow Fill with Dark Yellow Text" - Entre 30.000 e 120.000 (Regular)
$fcMedio = $estoque.FormatConditions.Add(1, 1, 30000, 120000)
$fcMedio.Font.Color = 22428
$fcMedio.Interior.Color = 10284031

# --- Restore selection ---------------------------------------------------
$ws.Range("C9").Select()
